$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: siddharth / siddhu / "1" / "1"
$ws.Range("A2").Value = "siddharth"
$ws.Range("B2").Value = "siddhu"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1"
$ws.Range("D2").Style = "Normal"

# Row 3 (new): nitesh / vada pav / "2" / "2"
$ws.Range("A3").Value = "nitesh"
$ws.Range("B3").Value = "vada pav"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2"
$ws.Range("D3").Style = "Normal"
